# CMPG323 EcoPower Logistics Data.xlsx
# "Added Choise to select mechanism"
#
# Fills in the boolean "Test Result" column on the Customers and Products
# sheets, then leaves the Products selection on G19 and the Customers tab
# as the active/selected one (matching the author's last actions before
# saving).

$wb = $excel.ActiveWorkbook

# --- Products sheet: mark the first 11 data rows (rows 2-12) as TRUE in
#     the "Test Result" column (E), then move the selection to G19.
$wsProducts = $wb.Worksheets.Item("Products")
$wsProducts.Activate()
$wsProducts.Range("E2:E12").Value = $true
$wsProducts.Range("G19").Select()

# --- Customers sheet: mark the first 5 data rows (rows 2-6) as TRUE in
#     the "Test Result" column (F). Activating this sheet last makes it
#     the one that is tabSelected / active when the workbook is saved.
$wsCustomers = $wb.Worksheets.Item("Customers")
$wsCustomers.Activate()
$wsCustomers.Range("F2:F6").Value = $true
